$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PeriodeItemStok_listItem")
$ws.Columns.Item(3).Insert()
$ws.Range("C1").ColumnWidth = 10.140625
Write-Host "done"
